$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 12, shifting the existing rows 12-58 down to 13-59.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new part entry (Retainer nut).
$ws.Range("A12").Value = "AFL-101-002-01-4"
$ws.Range("B12").Value = "A"
$ws.Range("C12").Value = "SLA in Rigid10k"
$ws.Range("E12").Value = "LoaderV2 (101)"
$ws.Range("F12").Value = "Catch Carrier Assy (002)"
$ws.Range("G12").Value = "Retainer nut"
$ws.Range("I12").Value = "Tap hole 4-40"

# Column I carries a workbook-level default style (no explicit per-cell
# style in the source data) - reset back to Normal so the newly written
# cell doesn't pick up an explicit style index from the row-insert.
$ws.Range("I12").Style = "Normal"

# Match the author's final cursor position.
$ws.Range("I8").Select()
